$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.197.71'
$ws.Range('E2').Value = '  -0.75%  '
$ws.Range('D3').Value = '2.595.54'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''310.03'
$ws.Range('E5').Value = '  +1.91%  '
$ws.Range('D6').Value = '''99.07'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  -0.45%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '''0.580'
$ws.Range('E9').Value = '  +0.63%  '
$ws.Range('D10').Value = '''38.96'
$ws.Range('E10').Value = '  +1.27%  '
$ws.Range('D11').Value = '''54.16'
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').Value = '''0.0839'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '2.991.48'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '2.599.63'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('E17').Value = '  +1.58%  '
$ws.Range('D18').Value = '''14.86'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '46.274.00'
$ws.Range('E19').Value = '  -0.77%  '
$ws.Range('E20').Value = '  +0.64%  '
$ws.Range('D21').Value = '''12.85'
$ws.Range('E21').Value = '  -5.06%  '
$ws.Range('D22').Value = '''6.71'
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').Value = '''71.55'
$ws.Range('E23').Value = '  +2.08%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '''272.94'
$ws.Range('E24').Value = '  +7.95%  '
$ws.Range('D25').Value = '''3.04'
$ws.Range('E25').Value = '  +2.89%  '
$ws.Range('E26').Value = '  +2.52%  '
$ws.Range('D27').Value = '''29.50'
$ws.Range('E27').Value = '  +9.26%  '
$ws.Range('D28').Value = '''0.998'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = '''4.07'
$ws.Range('E29').Value = '  +2.09%  '
$ws.Range('D30').Value = '''10.82'
$ws.Range('E30').Value = '  +3.96%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').Value = '''38.18'
$ws.Range('E31').Value = '  -2.55%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').Value = '''2.21'
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('D33').Value = '''6.26'
$ws.Range('E33').Value = '  +3.89%  '
$ws.Range('D34').Value = '''3.58'
$ws.Range('E34').Value = '  -6.98%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').Value = '''155.45'
$ws.Range('E35').Value = '  +3.41%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '''2.22'
$ws.Range('E36').Value = '  -4.56%  '
$ws.Range('D37').Value = '''0.0835'
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').Value = '''2.79'
$ws.Range('E38').Value = '  -5.75%  '
$ws.Range('D39').Value = '''0.123'
$ws.Range('E39').Value = '  +4.37%  '
$ws.Range('E40').Value = '  +1.09%  '
$ws.Range('D41').Value = '''22.43'
$ws.Range('E41').Value = '  +24.32%  '
$ws.Range('E42').Value = '  -0.05%  '
$ws.Range('E43').Value = '  +2.01%  '
$ws.Range('D44').Value = '''3.58'
$ws.Range('E44').Value = '  -0.95%  '
$ws.Range('D45').Value = '''3.96'
$ws.Range('E45').Value = '  -5.09%  '
$ws.Range('D46').Value = '2.109.27'
$ws.Range('E46').Value = '  +4.51%  '
$ws.Range('D47').Value = '''1.00'
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('D48').Value = '''94.98'
$ws.Range('E48').Value = '  +4.44%  '
$ws.Range('D49').Value = '''9.61'
$ws.Range('E49').Value = '  +7.33%  '
$ws.Range('D50').Value = '''108.68'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').Value = '''1.75'
$ws.Range('E51').Value = '  -1.80%  '
